# This script applies the TPM-based recalculation of ligand/receptor
# expression statistics (columns E-T) in the NATMI cell-cell signaling
# output sheet. Column values are updated to their new data-driven values;
# the edge weight/specificity columns (Q-T) are derivative of the
# ligand (G/H/I/J) and receptor (M/N/O/P) columns (Q=G*M, R=H*N, S=I*O, T=J*P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 0.3802213333333334
$ws.Range("H2").Value2 = 1.140664
$ws.Range("I2").Value2 = 0.05154022338265814
$ws.Range("J2").Value2 = 0.05154022338265814
$ws.Range("M2").Value2 = 10.055569
$ws.Range("N2").Value2 = 30.166707
$ws.Range("O2").Value2 = 0.349442268297237
$ws.Range("P2").Value2 = 0.3494422682972371
$ws.Range("Q2").Value2 = 3.823341852605334
$ws.Range("R2").Value2 = 34.410076673448
$ws.Range("S2").Value2 = 0.01801033256738236
$ws.Range("T2").Value2 = 0.01801033256738236
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 0.3802213333333334
$ws.Range("H3").Value2 = 1.140664
$ws.Range("I3").Value2 = 0.05154022338265814
$ws.Range("J3").Value2 = 0.05154022338265814
$ws.Range("O3").Value2 = 0.2478120087748427
$ws.Range("P3").Value2 = 0.2478120087748427
$ws.Range("Q3").Value2 = 2.711377846028444
$ws.Range("R3").Value2 = 24.402400614256
$ws.Range("S3").Value2 = 0.01277228628916063
$ws.Range("T3").Value2 = 0.01277228628916063
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 0.3802213333333334
$ws.Range("H4").Value2 = 1.140664
$ws.Range("I4").Value2 = 0.05154022338265814
$ws.Range("J4").Value2 = 0.05154022338265814
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 2.623176333333333
$ws.Range("N4").Value2 = 7.869529
$ws.Range("O4").Value2 = 0.09115831118692827
$ws.Range("P4").Value2 = 0.09115831118692828
$ws.Range("Q4").Value2 = 0.9973876030284445
$ws.Range("R4").Value2 = 8.976488427256001
$ws.Range("S4").Value2 = 0.004698319721760148
$ws.Range("T4").Value2 = 0.004698319721760148
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 0.3802213333333334
$ws.Range("H5").Value2 = 1.140664
$ws.Range("I5").Value2 = 0.05154022338265814
$ws.Range("J5").Value2 = 0.05154022338265814
$ws.Range("M5").Value2 = 6.063478666666666
$ws.Range("N5").Value2 = 18.190436
$ws.Range("O5").Value2 = 0.2107126646987263
$ws.Range("P5").Value2 = 0.2107126646987263
$ws.Range("Q5").Value2 = 2.305463943278222
$ws.Range("R5").Value2 = 20.749175489504
$ws.Range("S5").Value2 = 0.0108601778081275
$ws.Range("T5").Value2 = 0.0108601778081275
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 0.3802213333333334
$ws.Range("H6").Value2 = 1.140664
$ws.Range("I6").Value2 = 0.05154022338265814
$ws.Range("J6").Value2 = 0.05154022338265814
$ws.Range("K6").Value2 = 1
$ws.Range("L6").Value2 = 0.3333333333333333
$ws.Range("M6").Value2 = 0.5298106666666667
$ws.Range("N6").Value2 = 1.589432
$ws.Range("O6").Value2 = 0.01841151317524362
$ws.Range("P6").Value2 = 0.01841151317524363
$ws.Range("Q6").Value2 = 0.2014453180942222
$ws.Range("R6").Value2 = 1.813007862848
$ws.Range("S6").Value2 = 0.0009489335018648098
$ws.Range("T6").Value2 = 0.0009489335018648099
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 0.3802213333333334
$ws.Range("H7").Value2 = 1.140664
$ws.Range("I7").Value2 = 0.05154022338265814
$ws.Range("J7").Value2 = 0.05154022338265814
$ws.Range("M7").Value2 = 2.372966333333333
$ws.Range("N7").Value2 = 7.118898999999999
$ws.Range("O7").Value2 = 0.08246323386702208
$ws.Range("P7").Value2 = 0.08246323386702209
$ws.Range("Q7").Value2 = 0.902252423215111
$ws.Range("R7").Value2 = 8.120271808936
$ws.Range("S7").Value2 = 0.004250173494362698
$ws.Range("T7").Value2 = 0.004250173494362699
$ws.Range("G8").Value2 = 6.845175999999999
$ws.Range("H8").Value2 = 20.535528
$ws.Range("I8").Value2 = 0.9278856003177367
$ws.Range("J8").Value2 = 0.9278856003177366
$ws.Range("M8").Value2 = 10.055569
$ws.Range("N8").Value2 = 30.166707
$ws.Range("O8").Value2 = 0.349442268297237
$ws.Range("P8").Value2 = 0.3494422682972371
$ws.Range("Q8").Value2 = 68.83213958514399
$ws.Range("R8").Value2 = 619.4892562662959
$ws.Range("S8").Value2 = 0.3242424488953734
$ws.Range("T8").Value2 = 0.3242424488953734
$ws.Range("G9").Value2 = 6.845175999999999
$ws.Range("H9").Value2 = 20.535528
$ws.Range("I9").Value2 = 0.9278856003177367
$ws.Range("J9").Value2 = 0.9278856003177366
$ws.Range("O9").Value2 = 0.2478120087748427
$ws.Range("P9").Value2 = 0.2478120087748427
$ws.Range("Q9").Value2 = 48.81330144170133
$ws.Range("R9").Value2 = 439.319712975312
$ws.Range("S9").Value2 = 0.2299411945279892
$ws.Range("T9").Value2 = 0.2299411945279892
$ws.Range("G10").Value2 = 6.845175999999999
$ws.Range("H10").Value2 = 20.535528
$ws.Range("I10").Value2 = 0.9278856003177367
$ws.Range("J10").Value2 = 0.9278856003177366
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 2.623176333333333
$ws.Range("N10").Value2 = 7.869529
$ws.Range("O10").Value2 = 0.09115831118692827
$ws.Range("P10").Value2 = 0.09115831118692828
$ws.Range("Q10").Value2 = 17.95610368070133
$ws.Range("R10").Value2 = 161.604933126312
$ws.Range("S10").Value2 = 0.08458448429963399
$ws.Range("T10").Value2 = 0.084584484299634
$ws.Range("G11").Value2 = 6.845175999999999
$ws.Range("H11").Value2 = 20.535528
$ws.Range("I11").Value2 = 0.9278856003177367
$ws.Range("J11").Value2 = 0.9278856003177366
$ws.Range("M11").Value2 = 6.063478666666666
$ws.Range("N11").Value2 = 18.190436
$ws.Range("O11").Value2 = 0.2107126646987263
$ws.Range("P11").Value2 = 0.2107126646987263
$ws.Range("Q11").Value2 = 41.50557864557866
$ws.Range("R11").Value2 = 373.550207810208
$ws.Range("S11").Value2 = 0.1955172473785276
$ws.Range("T11").Value2 = 0.1955172473785276
$ws.Range("G12").Value2 = 6.845175999999999
$ws.Range("H12").Value2 = 20.535528
$ws.Range("I12").Value2 = 0.9278856003177367
$ws.Range("J12").Value2 = 0.9278856003177366
$ws.Range("K12").Value2 = 1
$ws.Range("L12").Value2 = 0.3333333333333333
$ws.Range("M12").Value2 = 0.5298106666666667
$ws.Range("N12").Value2 = 1.589432
$ws.Range("O12").Value2 = 0.01841151317524362
$ws.Range("P12").Value2 = 0.01841151317524363
$ws.Range("Q12").Value2 = 3.626647260010667
$ws.Range("R12").Value2 = 32.639825340096
$ws.Range("S12").Value2 = 0.01708377795536885
$ws.Range("T12").Value2 = 0.01708377795536885
$ws.Range("G13").Value2 = 6.845175999999999
$ws.Range("H13").Value2 = 20.535528
$ws.Range("I13").Value2 = 0.9278856003177367
$ws.Range("J13").Value2 = 0.9278856003177366
$ws.Range("M13").Value2 = 2.372966333333333
$ws.Range("N13").Value2 = 7.118898999999999
$ws.Range("O13").Value2 = 0.08246323386702208
$ws.Range("P13").Value2 = 0.08246323386702209
$ws.Range("Q13").Value2 = 16.24337219374133
$ws.Range("R13").Value2 = 146.190349743672
$ws.Range("S13").Value2 = 0.0765164472608437
$ws.Range("T13").Value2 = 0.0765164472608437
$ws.Range("G14").Value2 = 0.1517793333333333
$ws.Range("H14").Value2 = 0.455338
$ws.Range("I14").Value2 = 0.02057417629960514
$ws.Range("J14").Value2 = 0.02057417629960513
$ws.Range("M14").Value2 = 10.055569
$ws.Range("N14").Value2 = 30.166707
$ws.Range("O14").Value2 = 0.349442268297237
$ws.Range("P14").Value2 = 0.3494422682972371
$ws.Range("Q14").Value2 = 1.526227559107334
$ws.Range("R14").Value2 = 13.736048031966
$ws.Range("S14").Value2 = 0.007189486834481272
$ws.Range("T14").Value2 = 0.007189486834481272
$ws.Range("G15").Value2 = 0.1517793333333333
$ws.Range("H15").Value2 = 0.455338
$ws.Range("I15").Value2 = 0.02057417629960514
$ws.Range("J15").Value2 = 0.02057417629960513
$ws.Range("O15").Value2 = 0.2478120087748427
$ws.Range("P15").Value2 = 0.2478120087748427
$ws.Range("Q15").Value2 = 1.082346217339111
$ws.Range("R15").Value2 = 9.741115956052001
$ws.Range("S15").Value2 = 0.005098527957692909
$ws.Range("T15").Value2 = 0.005098527957692909
$ws.Range("G16").Value2 = 0.1517793333333333
$ws.Range("H16").Value2 = 0.455338
$ws.Range("I16").Value2 = 0.02057417629960514
$ws.Range("J16").Value2 = 0.02057417629960513
$ws.Range("K16").Value2 = 3
$ws.Range("L16").Value2 = 1
$ws.Range("M16").Value2 = 2.623176333333333
$ws.Range("N16").Value2 = 7.869529
$ws.Range("O16").Value2 = 0.09115831118692827
$ws.Range("P16").Value2 = 0.09115831118692828
$ws.Range("Q16").Value2 = 0.3981439550891112
$ws.Range("R16").Value2 = 3.583295595802
$ws.Range("S16").Value2 = 0.001875507165534129
$ws.Range("T16").Value2 = 0.001875507165534129
$ws.Range("G17").Value2 = 0.1517793333333333
$ws.Range("H17").Value2 = 0.455338
$ws.Range("I17").Value2 = 0.02057417629960514
$ws.Range("J17").Value2 = 0.02057417629960513
$ws.Range("M17").Value2 = 6.063478666666666
$ws.Range("N17").Value2 = 18.190436
$ws.Range("O17").Value2 = 0.2107126646987263
$ws.Range("P17").Value2 = 0.2107126646987263
$ws.Range("Q17").Value2 = 0.9203107497075556
$ws.Range("R17").Value2 = 8.282796747368
$ws.Range("S17").Value2 = 0.004335239512071177
$ws.Range("T17").Value2 = 0.004335239512071177
$ws.Range("G18").Value2 = 0.1517793333333333
$ws.Range("H18").Value2 = 0.455338
$ws.Range("I18").Value2 = 0.02057417629960514
$ws.Range("J18").Value2 = 0.02057417629960513
$ws.Range("K18").Value2 = 1
$ws.Range("L18").Value2 = 0.3333333333333333
$ws.Range("M18").Value2 = 0.5298106666666667
$ws.Range("N18").Value2 = 1.589432
$ws.Range("O18").Value2 = 0.01841151317524362
$ws.Range("P18").Value2 = 0.01841151317524363
$ws.Range("Q18").Value2 = 0.08041430977955556
$ws.Range("R18").Value2 = 0.723728788016
$ws.Range("S18").Value2 = 0.000378801718009965
$ws.Range("T18").Value2 = 0.000378801718009965
$ws.Range("G19").Value2 = 0.1517793333333333
$ws.Range("H19").Value2 = 0.455338
$ws.Range("I19").Value2 = 0.02057417629960514
$ws.Range("J19").Value2 = 0.02057417629960513
$ws.Range("M19").Value2 = 2.372966333333333
$ws.Range("N19").Value2 = 7.118898999999999
$ws.Range("O19").Value2 = 0.08246323386702208
$ws.Range("P19").Value2 = 0.08246323386702209
$ws.Range("Q19").Value2 = 0.3601672480957777
$ws.Range("R19").Value2 = 3.241505232862
$ws.Range("S19").Value2 = 0.001696613111815681
$ws.Range("T19").Value2 = 0.001696613111815681
